# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New values for rows 2-5, columns B-G
# (values passed as strings then cast to double to preserve exact precision,
#  including ones that require scientific notation)
$ws.Cells.Item(2, 2).Value = [double]"0.1169995834814548"
$ws.Cells.Item(2, 3).Value = [double]"6.708468553440206e-05"
$ws.Cells.Item(2, 4).Value = [double]"3.223369029078222"
$ws.Cells.Item(2, 5).Value = [double]"0.5333859586016987"
$ws.Cells.Item(2, 6).Value = [double]"0"
$ws.Cells.Item(2, 7).Value = [double]"3.87382165584691"

$ws.Cells.Item(3, 2).Value = [double]"1.445647641019636"
$ws.Cells.Item(3, 3).Value = [double]"1.626987699542094"
$ws.Cells.Item(3, 4).Value = [double]"0.1496068669990043"
$ws.Cells.Item(3, 5).Value = [double]"0.5333859586016987"
$ws.Cells.Item(3, 6).Value = [double]"0"
$ws.Cells.Item(3, 7).Value = [double]"3.755628166162433"

$ws.Cells.Item(4, 2).Value = [double]"1.445647641019636"
$ws.Cells.Item(4, 3).Value = [double]"1.626987699542094"
$ws.Cells.Item(4, 4).Value = [double]"3.223369029078222"
$ws.Cells.Item(4, 5).Value = [double]"0.5333859586016987"
$ws.Cells.Item(4, 6).Value = [double]"0"
$ws.Cells.Item(4, 7).Value = [double]"6.82939032824165"

$ws.Cells.Item(5, 2).Value = [double]"1.445647641019636"
$ws.Cells.Item(5, 3).Value = [double]"1.626987699542094"
$ws.Cells.Item(5, 4).Value = [double]"0.1496068669990043"
$ws.Cells.Item(5, 5).Value = [double]"0.5333859586016987"
$ws.Cells.Item(5, 6).Value = [double]"0"
$ws.Cells.Item(5, 7).Value = [double]"3.755628166162433"
